# Weekly fruit/vegetable price update: insert a new latest-week row at
# the top of the data (row 61), pushing all existing data rows down by
# one. The new row carries the newest week's Jengibre price data for
# Mercado Mayorista Lo Valledor de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 61; this shifts rows 61:91 down to 62:92
# and extends the used range to A1:R92.
$ws.Range("A61").EntireRow.Insert()

# Populate the newly inserted row 61 with the latest week's data.
$ws.Range("A61").Value = 6
$ws.Range("B61").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C61").Value = "Metropolitana"
$ws.Range("D61").Value = 44767
$ws.Range("E61").Value = 13
$ws.Range("F61").Value = 100114007
$ws.Range("G61").Value = "Jengibre"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 200
$ws.Range("K61").Value = 12000
$ws.Range("L61").Value = 13000
$ws.Range("M61").Value = 12600
$ws.Range("N61").Value = "`$/caja 13 kilos"
$ws.Range("O61").Value = "Per$([char]0x00FA)"
$ws.Range("P61").Value = 969
$ws.Range("Q61").Value = 13
$ws.Range("R61").Value = "Hortaliza"
